# Applies the "Artfynd" sheet update:
#  - rows 15-18 are re-sorted (their data is cyclically permuted) with a
#    couple of incidental field tweaks (an added/removed empty "Kön" cell,
#    and the "Publik kommentar" note moving along with its row's data)
#  - three brand new observation rows (19-21) are appended
#  - the sheet's used-range grows from A1:AY18 to A1:AY21 automatically as
#    a consequence of writing into row 21

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be stored as TEXT (prevents Excel's automatic
# number/date inference from mangling numeric-looking or date-looking
# strings such as "10" or "2023-09-26").
function Set-Text($range, $value) {
    $range.NumberFormat = "@"
    $range.Value2 = $value
    $range.ClearFormats()
}

# ---------------------------------------------------------------------
# Rows 15-18: the four observations get re-ordered. Columns that hold the
# exact same value in every one of the four rows (C, K, N, S, T, U, V, W,
# Y, AA, AD, AE, AF, AG, AT, AW, AX, AY) are intentionally left untouched.
# ---------------------------------------------------------------------

# New row 15 <- old row 17 data (Grönpyrola / Pyrola chlorantha)
$ws.Range("A15").Value2 = 111837675
$ws.Range("B15").Value2 = 103288
$ws.Range("D15").Value2 = "LC"
$ws.Range("E15").Value2 = 221144
$ws.Range("F15").Value2 = "Grönpyrola"
$ws.Range("G15").Value2 = "Pyrola chlorantha"
$ws.Range("H15").Value2 = "Sw."
Set-Text $ws.Range("I15") "10"
$ws.Range("J15").Value2 = "plantor/tuvor"
$ws.Range("P15").Value2 = "Brotorp, Långsjön, Sm"
$ws.Range("Q15").Value2 = 575782
$ws.Range("R15").Value2 = 6404547
# this row now needs the empty "Kön" placeholder column too
Set-Text $ws.Range("L15") ""

# New row 16 <- old row 15 data (Zontaggsvamp / Hydnellum concrescens)
$ws.Range("A16").Value2 = 111837705
$ws.Range("B16").Value2 = 90662
$ws.Range("D16").Value2 = "LC"
$ws.Range("E16").Value2 = 4363
$ws.Range("F16").Value2 = "Zontaggsvamp"
$ws.Range("G16").Value2 = "Hydnellum concrescens"
$ws.Range("H16").Value2 = "(Pers.) Banker"
Set-Text $ws.Range("I16") "10"
$ws.Range("J16").Value2 = "fruktkroppar"
$ws.Range("P16").Value2 = "Brotorp, Långsjön, Sm"
$ws.Range("Q16").Value2 = 575795
$ws.Range("R16").Value2 = 6404519

# New row 17 <- old row 18 data (Koralltaggsvamp / Hericium coralloides)
$ws.Range("A17").Value2 = 111837758
$ws.Range("B17").Value2 = 90187
$ws.Range("D17").Value2 = "NT"
$ws.Range("E17").Value2 = 2014
$ws.Range("F17").Value2 = "Koralltaggsvamp"
$ws.Range("G17").Value2 = "Hericium coralloides"
$ws.Range("H17").Value2 = "(Scop.:Fr.) Pers."
Set-Text $ws.Range("I17") "6"
$ws.Range("J17").Value2 = "fruktkroppar"
$ws.Range("P17").Value2 = "Brotorp, hyggeskant, Sm"
$ws.Range("Q17").Value2 = 575674
$ws.Range("R17").Value2 = 6404513
# this row no longer needs the empty "Kön" placeholder column
$ws.Range("L17").ClearContents()
# the public-comment note travels with the row it belongs to
Set-Text $ws.Range("AC17") "På asplåga."

# New row 18 <- old row 16 data (Orange taggsvamp / Hydnellum aurantiacum)
$ws.Range("A18").Value2 = 111837741
$ws.Range("B18").Value2 = 90658
$ws.Range("D18").Value2 = "NT"
$ws.Range("E18").Value2 = 4361
$ws.Range("F18").Value2 = "Orange taggsvamp"
$ws.Range("G18").Value2 = "Hydnellum aurantiacum"
$ws.Range("H18").Value2 = "(Batsch:Fr.) P.Karst."
Set-Text $ws.Range("I18") "15"
$ws.Range("J18").Value2 = "fruktkroppar"
$ws.Range("P18").Value2 = "Brotorp, hyggeskant, Sm"
$ws.Range("Q18").Value2 = 575654
$ws.Range("R18").Value2 = 6404507
# the public-comment note no longer belongs on this row
$ws.Range("AC18").ClearContents()

# ---------------------------------------------------------------------
# New rows 19-21: three additional observations at "Brotorp, Hummelstad"
# reported on 2023-09-26.
# ---------------------------------------------------------------------

# Row 19: Svart taggsvamp / Phellodon niger
$ws.Range("A19").Value2 = 112360063
$ws.Range("B19").Value2 = 90844
$ws.Range("C19").Value2 = "Ovaliderad"
$ws.Range("D19").Value2 = "NT"
$ws.Range("E19").Value2 = 5449
$ws.Range("F19").Value2 = "Svart taggsvamp"
$ws.Range("G19").Value2 = "Phellodon niger"
$ws.Range("H19").Value2 = "(Fr.:Fr.) P.Karst."
Set-Text $ws.Range("I19") "2"
$ws.Range("J19").Value2 = "fruktkroppar"
$ws.Range("P19").Value2 = "Brotorp, Hummelstad, Sm"
$ws.Range("Q19").Value2 = 575771
$ws.Range("R19").Value2 = 6404602
$ws.Range("S19").Value2 = 10
$ws.Range("T19").Value2 = "Kalmar"
$ws.Range("U19").Value2 = "Västervik"
$ws.Range("V19").Value2 = "Småland"
$ws.Range("W19").Value2 = "Hallingeberg"
Set-Text $ws.Range("Y19") "2023-09-26"
Set-Text $ws.Range("AA19") "2023-09-26"
$ws.Range("AD19").Value2 = $false
$ws.Range("AE19").Value2 = $false
$ws.Range("AG19").Value2 = $false
$ws.Range("AW19").Value2 = "Magnus Kasselstrand"
$ws.Range("AX19").Value2 = "Magnus Kasselstrand"

# Row 20: Grönpyrola / Pyrola chlorantha (second location)
$ws.Range("A20").Value2 = 112360070
$ws.Range("B20").Value2 = 103755
$ws.Range("C20").Value2 = "Ovaliderad"
$ws.Range("D20").Value2 = "LC"
$ws.Range("E20").Value2 = 221144
$ws.Range("F20").Value2 = "Grönpyrola"
$ws.Range("G20").Value2 = "Pyrola chlorantha"
$ws.Range("H20").Value2 = "Sw."
$ws.Range("J20").Value2 = "plantor/tuvor"
$ws.Range("P20").Value2 = "Brotorp, Hummelstad, Sm"
$ws.Range("Q20").Value2 = 575771
$ws.Range("R20").Value2 = 6404602
$ws.Range("S20").Value2 = 10
$ws.Range("T20").Value2 = "Kalmar"
$ws.Range("U20").Value2 = "Västervik"
$ws.Range("V20").Value2 = "Småland"
$ws.Range("W20").Value2 = "Hallingeberg"
Set-Text $ws.Range("Y20") "2023-09-26"
Set-Text $ws.Range("AA20") "2023-09-26"
$ws.Range("AC20").Value2 = "Få"
$ws.Range("AD20").Value2 = $false
$ws.Range("AE20").Value2 = $false
$ws.Range("AG20").Value2 = $false
$ws.Range("AW20").Value2 = "Magnus Kasselstrand"
$ws.Range("AX20").Value2 = "Magnus Kasselstrand"

# Row 21: Dropptaggsvamp / Hydnellum ferrugineum
$ws.Range("A21").Value2 = 112360057
$ws.Range("B21").Value2 = 90800
$ws.Range("C21").Value2 = "Ovaliderad"
$ws.Range("D21").Value2 = "LC"
$ws.Range("E21").Value2 = 4364
$ws.Range("F21").Value2 = "Dropptaggsvamp"
$ws.Range("G21").Value2 = "Hydnellum ferrugineum"
$ws.Range("H21").Value2 = "(Fr.:Fr.) P. Karst."
Set-Text $ws.Range("I21") "27"
$ws.Range("J21").Value2 = "fruktkroppar"
$ws.Range("P21").Value2 = "Brotorp, Hummelstad, Sm"
$ws.Range("Q21").Value2 = 575771
$ws.Range("R21").Value2 = 6404602
$ws.Range("S21").Value2 = 10
$ws.Range("T21").Value2 = "Kalmar"
$ws.Range("U21").Value2 = "Västervik"
$ws.Range("V21").Value2 = "Småland"
$ws.Range("W21").Value2 = "Hallingeberg"
Set-Text $ws.Range("Y21") "2023-09-26"
Set-Text $ws.Range("AA21") "2023-09-26"
$ws.Range("AD21").Value2 = $false
$ws.Range("AE21").Value2 = $false
$ws.Range("AG21").Value2 = $false
$ws.Range("AW21").Value2 = "Magnus Kasselstrand"
$ws.Range("AX21").Value2 = "Magnus Kasselstrand"

Write-Output "edit complete"
